{"js": "// V15 - \"Strategist gebruiken is geen actie.\"\n// Rewrites the Strategist ability description: it becomes something the\n// Strategist does at the start of the turn (not an action that replaces a\n// move/attack), and \"onthullen\" (reveal) is reworded to \"bekijken\" (view).\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"Strategist\" Heading 2 paragraph, then work on the three\n// paragraphs that directly follow it (the ones describing its ability).\nlet stratHeadingIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const item = paras.items[i];\n  if (item.text.trim() === \"Strategist\" && item.style && item.style.indexOf(\"Heading\") !== -1) {\n    stratHeadingIdx = i;\n    break;\n  }\n}\nif (stratHeadingIdx === -1) {\n  throw new Error(\"Could not find the 'Strategist' heading paragraph.\");\n}\n\nconst p1 = paras.items[stratHeadingIdx + 1]; // \"Na zich getoond te hebben kan de Strategist ...\"\nconst p2 = paras.items[stratHeadingIdx + 2]; // \"Dit wordt gedaan in plaats van bewegen of slaan.\"\nconst p3 = paras.items[stratHeadingIdx + 3]; // \"Deze actie negeert rang. ...\"\n\n// --- Paragraph 1: full rewrite ---\np1.insertText(\n  \"Aan het begin van de beurt kan de Strategist zich tonen om een vijandelijk stuk dat het kan bereiken, te bekijken; de eigenaar moet het karakter van het stuk laten zien. \",\n  \"Replace\"\n);\n\n// --- Paragraph 2: full rewrite ---\np2.insertText(\n  \"Dit gebeurt voordat de speler een stuk gebruikt om te bewegen of te slaan. De Strategist mag maar \u00e9\u00e9n stuk per beurt bekijken.\",\n  \"Replace\"\n);\n\nawait context.sync();\n\n// --- Paragraph 3: move the _GoBack bookmark + reword \"onthullen\" -> \"bekijken\" ---\n\n// The _GoBack bookmark used to sit (empty) right before \" 2 velden ver\n// bevinden.\"; it now wraps the \"Deze actie negeert rang. \" run at the start\n// of the paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst rangeNegeertRang = p3.search(\"Deze actie negeert rang. \", { matchCase: true });\nrangeNegeertRang.load(\"items/text\");\nawait context.sync();\nif (rangeNegeertRang.items.length > 0) {\n  rangeNegeertRang.items[0].insertBookmark(\"_GoBack\");\n}\n\n// \"... stuk in de Tuin onthullen.\" -> \"... stuk in de Tuin bekijken.\"\nconst rangeTuin = p3.search(\"stuk in de Tuin onthullen\", { matchCase: true });\nrangeTuin.load(\"items/text\");\nawait context.sync();\nif (rangeTuin.items.length > 0) {\n  rangeTuin.items[0].insertText(\"stuk in de Tuin \", \"Replace\");\n  await context.sync();\n  const rangeTuinTail = p3.search(\"stuk in de Tuin \", { matchCase: true });\n  rangeTuinTail.load(\"items/text\");\n  await context.sync();\n  rangeTuinTail.items[0].insertText(\"bekijken\", \"After\");\n}\n\n// \"... geen stukken onthullen die zich meer dan 2 ...\" -> \"... bekijken die zich ...\"\nconst rangeStukkenOnthullen = p3.search(\"onthullen die zich\", { matchCase: true });\nrangeStukkenOnthullen.load(\"items/text\");\nawait context.sync();\nif (rangeStukkenOnthullen.items.length > 0) {\n  rangeStukkenOnthullen.items[0].insertText(\"bekijken die zich\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# V15 - \"Strategist gebruiken is geen actie.\"\n# Rewrites the Strategist ability description: it becomes something the\n# Strategist does at the start of the turn (not an action that replaces a\n# move/attack), and \"onthullen\" (reveal) is reworded to \"bekijken\" (view).\n\n$d = $word.ActiveDocument\n\n# Locate the \"Strategist\" Heading 2 paragraph, then work on the three\n# paragraphs that directly follow it (the ones describing its ability).\n$headingIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $txt = $para.Range.Text.Trim()\n    if ($txt -eq \"Strategist\" -and $para.Style.NameLocal -like \"*Heading*\") {\n        $headingIdx = $i\n        break\n    }\n}\nif ($headingIdx -eq -1) {\n    throw \"Could not find the 'Strategist' heading paragraph.\"\n}\n\n$p1 = $d.Paragraphs.Item($headingIdx + 1).Range   # \"Na zich getoond te hebben kan de Strategist ...\"\n$p2 = $d.Paragraphs.Item($headingIdx + 2).Range   # \"Dit wordt gedaan in plaats van bewegen of slaan.\"\n$p3 = $d.Paragraphs.Item($headingIdx + 3).Range   # \"Deze actie negeert rang. ...\"\n\n# --- Paragraph 1: full rewrite ---\n$find1 = $p1.Find\n$find1.ClearFormatting()\n$find1.Execute(\n    \"Na zich getoond te hebben kan de Strategist een vijandelijk stuk dat het kan bereiken onthullen; de eigenaar moet het karakter van het stuk laten zien.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Aan het begin van de beurt kan de Strategist zich tonen om een vijandelijk stuk dat het kan bereiken, te bekijken; de eigenaar moet het karakter van het stuk laten zien. \",\n    2\n) | Out-Null\n\n# --- Paragraph 2: full rewrite ---\n$find2 = $p2.Find\n$find2.ClearFormatting()\n$find2.Execute(\n    \"Dit wordt gedaan in plaats van bewegen of slaan.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Dit gebeurt voordat de speler een stuk gebruikt om te bewegen of te slaan. De Strategist mag maar \u00e9\u00e9n stuk per beurt bekijken.\",\n    2\n) | Out-Null\n\n# --- Paragraph 3: reword \"onthullen\" -> \"bekijken\" (twice) and move the _GoBack bookmark ---\n\n# \"... stuk in de Tuin onthullen.\" -> \"... stuk in de Tuin bekijken.\"\n# \"... geen stukken onthullen die zich ...\" -> \"... geen stukken bekijken die zich ...\"\n$find3 = $p3.Find\n$find3.ClearFormatting()\n$find3.Execute(\n    \"onthullen\", $true, $false, $false, $false, $false, $true, 1, $false,\n    \"bekijken\", 2\n) | Out-Null\n\n# The _GoBack bookmark used to sit (empty) right before \" 2 velden ver\n# bevinden.\"; it now wraps the \"Deze actie negeert rang. \" run at the start\n# of the paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$bookmarkRange = $p3.Duplicate\n$bookmarkRange.Find.Execute(\"Deze actie negeert rang. \") | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n"}
